{"js": "// Find the list item that currently reads:\n//   \"Cutting more distance than others of the same day =+ 1p\"\n// and replace its text with:\n//   \"Most rides in a week =+ 2p\"\n// The run formatting (font, size, color) is preserved automatically because\n// Range.insertText(\"Replace\") keeps the existing run properties of the\n// paragraph it rewrites.\n\nconst oldText = \"Cutting more distance than others of the same day =+ 1p\";\nconst newText = \"Most rides in a week =+ 2p\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find paragraph text: \"${oldText}\"`);\n}\n\n// Replace the whole paragraph's range (not just the matched search hit) so\n// any trailing content/marks stay intact and the run keeps its original\n// formatting.\nconst hit = results.items[0];\nconst paragraph = hit.paragraphs.getFirst();\nconst fullRange = paragraph.getRange();\nfullRange.insertText(newText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Replace the ranking-criteria bullet:\n#   \"Cutting more distance than others of the same day =+ 1p\"\n# with:\n#   \"Most rides in a week =+ 2p\"\n# Find & Replace keeps the existing run formatting (font, size, color) of the\n# paragraph being rewritten.\n\n$d = $word.ActiveDocument\n\n$oldText = \"Cutting more distance than others of the same day =+ 1p\"\n$newText = \"Most rides in a week =+ 2p\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# Wrap=0 (wdFindStop), Replace=2 (wdReplaceAll)\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n"}
